$wb = $excel.ActiveWorkbook

# --- "beads" sheet: insert a new column C for "FL3-H Peaks" ---
$wsBeads = $wb.Worksheets.Item("beads")
$wsBeads.Columns("C").Insert()
$wsBeads.Range("C1").Value = "FL3-H Peaks"
$wsBeads.Range("C2").Value = $wsBeads.Range("B2").Value2
$wsBeads.Columns("C").ColumnWidth = $wsBeads.Columns("B").ColumnWidth
$wsBeads.Range("C2").Select()

# --- "cells" sheet: insert a new column E for "FL3-H Transform" ---
$wsCells = $wb.Worksheets.Item("cells")
$wsCells.Columns("E").Insert()
$wsCells.Range("E1").Value = "FL3-H Transform"
$wsCells.Range("E6").Value = $wsCells.Range("C6").Value2
$wsCells.Range("E1").Select()
